# Update cfb_weather.xlsx with Timestamp 2024-12-16T10:01:31.821318
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FBS")

# 1) Refresh the shared "Timestamp" string used across the whole AK column.
#    (All AK cells reference the same shared string, so updating the
#    column in one pass keeps the workbook's shared string table correct.)
$oldTimestamp = "2024-12-16T08:54:06.481427"
$newTimestamp = "2024-12-16T10:01:31.821318"

$used = $ws.UsedRange
$lastRow = $used.Rows.Count
for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Range("AK" + $r)
    if ($cell.Value() -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}

# 2) Row 8 weather values
$ws.Range("M8").Value = "N"
$ws.Range("N8").Value = "NNW"
$ws.Range("O8").Value = 37.76
$ws.Range("P8").Value = 3.1
$ws.Range("Q8").Value = "NNW"
$ws.Range("U8").Value = -2.6

# 3) Row 15 weather values
$ws.Range("AB15").Value = 7.5
$ws.Range("AF15").Value = 0.5

# 4) Row 16 weather values
$ws.Range("Z16").Value = -105

# 5) Row 18 weather values
$ws.Range("AB18").Value = 8.5
$ws.Range("AF18").Value = -1
